# Removed the old way of downloading full albums (dynamically created zip
# files) from the backlog. This deletes the whole "Zip files" row (row 3:
# "Zip files" / "Use real zip files rather than ones created on the fly or
# make sure these don't fail to download") and shifts all following rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 3 - this shifts rows 4..16 up to 3..15
$ws.Rows.Item(3).Delete()

# Leave selection on the sheet as Excel would after the delete operation
$ws.Range("B22").Select()
